# Updates to pl_mw.xlsx results sheet ("case with 380 kV done")
# Rewrites columns B, D, E, F, G, I, K, L for rows 2-25 with new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 4, 5, 6, 7, 9, 11, 12)   # B, D, E, F, G, I, K, L
$startRow = 2

$data = @(
    @(0.5190871431972823,0.1190062623558674,0.06415857793147062,2.466724915107989,0.002532243446708527,1.824092820839297,0.8231104967113367,0.2950475121658513),
    @(0.5087830456629234,0.118189322361637,0.06298639832995256,2.399267520623908,0.002537325230754242,1.800343700363314,0.7540332428670524,0.2814279673159206),
    @(0.502835730117468,0.1176923251201174,0.06225217663774796,2.358830509865385,0.002540609031785225,1.786195376329331,0.7121914339318209,0.273256141499445),
    @(0.5005077799678048,0.1174909165794986,0.06194926451315119,2.342597567824996,0.002541988477594099,1.780538335155541,0.695283139665861,0.2699738035559847),
    @(0.5001270078269755,0.1174575395289281,0.06189874017349872,2.339916882555968,0.002542220030228547,1.779605525830796,0.6924841065041107,0.2694316522549229),
    @(0.5028039470807215,0.117689604360919,0.06224810656772028,2.358610594825691,0.002540627468184863,1.786118644857268,0.7119628270573628,0.2732116816603423),
    @(0.5154556598492093,0.1187236029720822,0.06375738006666509,2.443260847141062,0.00253396178274967,1.815813763245146,0.7991731937729014,0.2903118089038088),
    @(0.5432692788095039,0.1207895624821447,0.06660505690919472,2.617129846088716,0.002522181797816957,1.877515772940654,0.9747959200570335,0.3253683676232555),
    @(0.5655291410431289,0.1223332863957936,0.0686333193248263,2.749791626285372,0.002514305326397243,1.925008446808476,1.106738105527029,0.3520710614844518),
    @(0.5760511131559554,0.1230417302357978,0.06954310741844871,2.811238749886883,0.002510889192652278,1.947093161728063,1.167418412017696,0.3644284958955097),
    @(0.5800922929035153,0.1233109318220471,0.06988584507777951,2.834667088191736,0.002509619446701548,1.955525798044846,1.190492933551184,0.3691384431745632),
    @(0.5792194320292481,0.1232529125032755,0.06981210836023699,2.829614249310026,0.002509891849672917,1.953706572314701,1.185519120734796,0.3681227136763709),
    @(0.5763824473901877,0.1230638587688375,0.06957133989288344,2.813163004223611,0.002510784252101152,1.947785520762721,1.169314830426231,0.3648153744401839),
    @(0.5746520958371661,0.1229481800524113,0.06942363270443508,2.803106981997388,0.002511333979076991,1.944167792055566,1.159401806321227,0.3627935064743326),
    @(0.5648494549523377,0.122287115879061,0.06857360956203173,2.745798161259415,0.002514531926176074,1.923574870171095,1.1027858894102,0.3512677236360844),
    @(0.5589371051873115,0.1218831886284377,0.06804890066971936,2.710923617907753,0.002516536419444523,1.911065161753911,1.068223684126508,0.3442510346197594),
    @(0.5555737479335789,0.1216514412739826,0.06774588538564785,2.690968009077437,0.002517705069326997,1.903915057633412,1.048406372863752,0.3402349963984932),
    @(0.5544413798386358,0.1215730742087899,0.06764307847037188,2.684229075220628,0.002518103457816679,1.901501891365129,1.041707169738572,0.3388786246634652),
    @(0.5595626289747599,0.1219261269738467,0.06810488221895028,2.714625368777433,0.002516321411793206,1.912392164440092,1.071896467551881,0.3449959248380026),
    @(0.5772142001663951,0.1231193629076728,0.06964210713521446,2.817990787566345,0.00251052148516385,1.94952278246781,1.174071800801357,0.3657859918118618),
    @(0.5890811482526601,0.1239046458992377,0.07063643556256238,2.886477390333312,0.002506869969790396,1.974195950244635,1.241410732084489,0.3795511240159897),
    @(0.5827173423069496,0.1234850156924097,0.0701066653270157,2.84983902976748,0.002508806170404052,1.960990052561002,1.205418827769449,0.3721880869854033),
    @(0.5592797183608411,0.1219067130689311,0.06807957718075563,2.71295151382725,0.002516418566135238,1.911792096080987,1.070235837310804,0.3446591040775075),
    @(0.5354240267756438,0.120226305766252,0.06584622919717198,2.569238531678337,0.002525231273894821,1.860447636629829,0.9267816492443899,0.3157197166221835)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowVals[$j]
    }
}
